$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "MSAT30"
$ws.Range("B2").Value = "Au.Panamericana _Blas Parera 1297|"
$ws.Range("C2").Value = -34.5354354
$ws.Range("D2").Value = -58.5021847
$ws.Range("E2").Value = "https://i.ibb.co/x215Lvb/Blas-Parera-1267-Au-Panamericana.jpg"
$ws.Range("F2").Value = "Vicente Lopez"
$ws.Range("H2").Value = "BUENOS AIRES"
$ws.Range("I2").Value = "GBA NORTE"
$ws.Range("L2").Value = "Cartel Espectacular Doble Faz"

# Row 3
$ws.Range("A3").Value = "MSAT33"
$ws.Range("B3").Value = "Av. Crisólogo Larralde 899"
$ws.Range("C3").Value = -34.6766898
$ws.Range("D3").Value = -58.4598458
$ws.Range("E3").Value = "https://i.ibb.co/q0v38Mj/Av-Cris-logo-Larralde-899-Avellaneda-GBA.jpg"
$ws.Range("F3").Value = "Avellaneda"
$ws.Range("H3").Value = "BUENOS AIRES"
$ws.Range("I3").Value = "GBA SUR"
$ws.Range("L3").Value = "Cartel Espectacular Simple Faz"

# Row 4
$ws.Range("A4").Value = "MSAT34"
$ws.Range("B4").Value = "Av. Del Libertador 240"
$ws.Range("C4").Value = -34.532315
$ws.Range("D4").Value = -58.4712099
$ws.Range("E4").Value = "https://i.ibb.co/VCHZWdj/Av-Del-Libertador-240-Vte-L-pez.jpg"
$ws.Range("F4").Value = "Vicente Lopez"
$ws.Range("H4").Value = "BUENOS AIRES"
$ws.Range("I4").Value = "GBA NORTE"

# Remove rows 5 and 6 (old MSAT46 / MA63 entries)
$ws.Rows("5:6").Delete()
